$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.774.48"
$ws.Range("E2").Value = "  -1.57%  "
$ws.Range("D3").Value = "2.234.67"
$ws.Range("E3").Value = "  -2.26%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "115.09"
$ws.Range("E5").Value = "  +1.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "266.93"
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("E7").Value = "  +1.17%  "
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.605"
$ws.Range("E9").Value = "  -0.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "46.42"
$ws.Range("E10").Value = "  -3.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0929"
$ws.Range("E11").Value = "  -0.42%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.26"
$ws.Range("E12").Value = "  -0.97%  "
$ws.Range("E13").Value = "  -2.83%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.35"
$ws.Range("E14").Value = "  -1.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.880"
$ws.Range("E15").Value = "  +1.32%  "
$ws.Range("D16").Value = "2.569.33"
$ws.Range("E16").Value = "  -1.89%  "
$ws.Range("D17").Value = "2.240.97"
$ws.Range("E17").Value = "  -2.53%  "
$ws.Range("D18").Value = "43.022.59"
$ws.Range("E18").Value = "  -0.82%  "
$ws.Range("E19").Value = "  -0.90%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.73"
$ws.Range("E20").Value = "  -1.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.60"
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("E22").Value = "  -5.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "231.41"
$ws.Range("E23").Value = "  -0.89%  "
$ws.Range("E24").Value = "  -1.52%  "
$ws.Range("E25").Value = "  -4.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.44"
$ws.Range("E26").Value = "  +8.57%  "
$ws.Range("E27").Value = "  +0.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "40.72"
$ws.Range("E28").Value = "  -0.48%  "
$ws.Range("E29").Value = "  -1.11%  "
$ws.Range("E30").Value = "  -1.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.37"
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.17"
$ws.Range("E32").Value = "  -1.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0901"
$ws.Range("E33").Value = "  -0.82%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.59"
$ws.Range("E34").Value = "  -2.76%  "
$ws.Range("E35").Value = "  +7.55%  "
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0373"
$ws.Range("E37").Value = "  +1.78%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.62"
$ws.Range("E38").Value = "  -1.00%  "
$ws.Range("E39").Value = "  +1.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.51"
$ws.Range("E40").Value = "  -6.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "13.40"
$ws.Range("E41").Value = "  -6.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "71.01"
$ws.Range("E42").Value = "  -9.57%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.233"
$ws.Range("E43").Value = "  -2.75%  "
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("E45").Value = "  -3.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.65"
$ws.Range("E46").Value = "  -10.16%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.655"
$ws.Range("E47").Value = "  +10.98%  "
$ws.Range("B48").Value = "TrustWalletToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.25"
$ws.Range("E48").Value = "  -0.92%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0992"
$ws.Range("E49").Value = "  -0.65%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.42"
$ws.Range("E50").Value = "  -3.36%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "100.26"
$ws.Range("E51").Value = "  -4.14%  "
